$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 18131.143
$ws.Range("I32").Value = 14119.857
$ws.Range("J32").Value = 22142.428
$ws.Range("K32").Value = 14119.857
$ws.Range("L32").Value = 22142.428
$ws.Range("M32").Value = -13793.857
$ws.Range("N32").Value = -22794.428

$ws.Range("H76").Value = 4425.619
$ws.Range("J76").Value = 5875.4443
$ws.Range("L76").Value = 5875.4443
$ws.Range("N76").Value = -6505.4443

$ws.Range("H79").Value = 4425.619
$ws.Range("J79").Value = 5875.4443
$ws.Range("L79").Value = 5875.4443
$ws.Range("N79").Value = -8059.4443

$ws.Range("H86").Value = 7349.65
$ws.Range("I86").Value = 7179.5
$ws.Range("J86").Value = 7519.8
$ws.Range("K86").Value = 7179.5
$ws.Range("L86").Value = 7519.8
$ws.Range("M86").Value = -6056.5
$ws.Range("N86").Value = -9765.799999999999

$ws.Range("H89").Value = 7349.65
$ws.Range("I89").Value = 7179.5
$ws.Range("J89").Value = 7519.8
$ws.Range("K89").Value = 35897.5
$ws.Range("L89").Value = 37599
$ws.Range("M89").Value = -30281.5
$ws.Range("N89").Value = -48831

$ws.Range("H98").Value = 76924810
$ws.Range("I98").Value = 83335080
$ws.Range("K98").Value = 83335080
$ws.Range("M98").Value = -83333582

$ws.Range("H122").Value = 76924810
$ws.Range("I122").Value = 83335080
$ws.Range("K122").Value = 250005240
$ws.Range("M122").Value = -250002790

$ws.Range("H129").Value = 1903.2
$ws.Range("I129").Value = 1512.5
$ws.Range("J129").Value = 2163.6667
$ws.Range("K129").Value = 4537.5
$ws.Range("L129").Value = 6491.000100000001
$ws.Range("M129").Value = 462.5
$ws.Range("N129").Value = -16491.0001

$ws.Range("H132").Value = 1211.2307
$ws.Range("I132").Value = 1211.2307
$ws.Range("K132").Value = 3633.6921
$ws.Range("M132").Value = -1103.6921

$ws.Range("H141").Value = 3019
$ws.Range("I141").Value = 3019
$ws.Range("K141").Value = 9057
$ws.Range("M141").Value = -3877

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11364584
$ws.Range("I32").Value = 11905535
$ws.Range("K32").Value = 11905535
$ws.Range("M32").Value = -11905248

$ws.Range("H45").Value = 1825.5454
$ws.Range("I45").Value = 796.1111
$ws.Range("K45").Value = 796.1111
$ws.Range("M45").Value = -419.1111

$ws.Range("H74").Value = 13169855
$ws.Range("I74").Value = 22729434
$ws.Range("J74").Value = 25433.25
$ws.Range("K74").Value = 22729434
$ws.Range("L74").Value = 25433.25
$ws.Range("M74").Value = -22728560
$ws.Range("N74").Value = -27181.25

$ws.Range("H77").Value = 13169855
$ws.Range("I77").Value = 22729434
$ws.Range("J77").Value = 25433.25
$ws.Range("K77").Value = 113647170
$ws.Range("L77").Value = 127166.25
$ws.Range("M77").Value = -113642802
$ws.Range("N77").Value = -135902.25

$ws.Range("H88").Value = 1536.6471
$ws.Range("I88").Value = 1329.4
$ws.Range("J88").Value = 1832.7142
$ws.Range("K88").Value = 1329.4
$ws.Range("L88").Value = 1832.7142
$ws.Range("M88").Value = -923.4000000000001
$ws.Range("N88").Value = -2644.7142

$ws.Range("H91").Value = 1536.6471
$ws.Range("I91").Value = 1329.4
$ws.Range("J91").Value = 1832.7142
$ws.Range("K91").Value = 1329.4
$ws.Range("L91").Value = 1832.7142
$ws.Range("M91").Value = 74.59999999999991
$ws.Range("N91").Value = -4640.7142

$ws.Range("H122").Value = 3875.2222
$ws.Range("I122").Value = 1626
$ws.Range("K122").Value = 4878
$ws.Range("M122").Value = -2428

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1765.9286
$ws.Range("I105").Value = 1526.6666
$ws.Range("J105").Value = 2196.6
$ws.Range("K105").Value = 1526.6666
$ws.Range("L105").Value = 2196.6
$ws.Range("M105").Value = 220.3334
$ws.Range("N105").Value = -5690.6

$ws.Range("H134").Value = 51643.74
$ws.Range("I134").Value = 6910.1577
$ws.Range("K134").Value = 20730.4731
$ws.Range("M134").Value = -18195.4731

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 3592
$ws.Range("I7").Value = 321
$ws.Range("J7").Value = 7517.2
$ws.Range("K7").Value = 321
$ws.Range("L7").Value = 7517.2
$ws.Range("M7").Value = -208
$ws.Range("N7").Value = -7743.2

$ws.Range("H16").Value = 2721.75
$ws.Range("I16").Value = 1962.6666
$ws.Range("J16").Value = 4999
$ws.Range("K16").Value = 1962.6666
$ws.Range("L16").Value = 4999
$ws.Range("M16").Value = -1675.6666
$ws.Range("N16").Value = -5573

$ws.Range("H22").Value = 939.3333
$ws.Range("I22").Value = 934.5
$ws.Range("J22").Value = 949
$ws.Range("K22").Value = 934.5
$ws.Range("L22").Value = 949
$ws.Range("M22").Value = -584.5
$ws.Range("N22").Value = -1649

$ws.Range("H31").Value = 1061084.5
$ws.Range("I31").Value = 1440.4667
$ws.Range("K31").Value = 1440.4667
$ws.Range("M31").Value = -1145.4667

$ws.Range("H34").Value = 1061084.5
$ws.Range("I34").Value = 1440.4667
$ws.Range("K34").Value = 1440.4667
$ws.Range("M34").Value = -1238.4667

$ws.Range("H113").Value = 2721.75
$ws.Range("I113").Value = 1962.6666
$ws.Range("J113").Value = 4999
$ws.Range("K113").Value = 1962.6666
$ws.Range("L113").Value = 4999
$ws.Range("M113").Value = 207.3334
$ws.Range("N113").Value = -9339

$ws.Range("H122").Value = 4859.7
$ws.Range("I122").Value = 1979.5714
$ws.Range("K122").Value = 5938.7142
$ws.Range("M122").Value = -3488.7142

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 10900036
$ws.Range("I4").Value = 9333378
$ws.Range("K4").Value = 28000134
$ws.Range("M4").Value = -28000022

$ws.Range("H11").Value = 412.5
$ws.Range("I11").Value = 300
$ws.Range("K11").Value = 900
$ws.Range("M11").Value = -760

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 90.583336
$ws.Range("J2").Value = 109.8
$ws.Range("L2").Value = 109.8
$ws.Range("N2").Value = -335.8

$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents()

$ws.Range("H70").Value = 4290.3335
$ws.Range("I70").Value = 3581.3333
$ws.Range("J70").Value = 4999.3335
$ws.Range("K70").Value = 3581.3333
$ws.Range("L70").Value = 4999.3335
$ws.Range("M70").Value = -3311.3333
$ws.Range("N70").Value = -5539.3335

$ws.Range("H73").Value = 4290.3335
$ws.Range("I73").Value = 3581.3333
$ws.Range("J73").Value = 4999.3335
$ws.Range("K73").Value = 3581.3333
$ws.Range("L73").Value = 4999.3335
$ws.Range("M73").Value = -2645.3333
$ws.Range("N73").Value = -6871.3335

$ws.Range("H122").Value = 3629.75
$ws.Range("I122").Value = 2955.8
$ws.Range("K122").Value = 8867.400000000001
$ws.Range("M122").Value = -6417.400000000001

$ws.Range("H123").Value = 39988
$ws.Range("J123").Value = 39988
$ws.Range("L123").Value = 39988
$ws.Range("N123").Value = -44888

$ws.Range("H133").Value = 51250
$ws.Range("J133").Value = 51250
$ws.Range("L133").Value = 51250
$ws.Range("N133").Value = -61370

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 50001.773
$ws.Range("J7").Value = 95889.73
$ws.Range("L7").Value = 95889.73
$ws.Range("N7").Value = -96113.73

$ws.Range("H22").Value = 2046.0667
$ws.Range("I22").Value = 2524.125
$ws.Range("J22").Value = 1499.7142
$ws.Range("K22").Value = 2524.125
$ws.Range("L22").Value = 1499.7142
$ws.Range("M22").Value = -2229.125
$ws.Range("N22").Value = -2089.7142

$ws.Range("H27").Value = 2046.0667
$ws.Range("I27").Value = 2524.125
$ws.Range("J27").Value = 1499.7142
$ws.Range("K27").Value = 2524.125
$ws.Range("L27").Value = 1499.7142
$ws.Range("M27").Value = -2417.125
$ws.Range("N27").Value = -1713.7142

$ws.Range("H68").Value = 4300.6665
$ws.Range("I68").Value = 1826
$ws.Range("J68").Value = 9250
$ws.Range("K68").Value = 1826
$ws.Range("L68").Value = 9250
$ws.Range("M68").Value = -1077
$ws.Range("N68").Value = -10748

$ws.Range("H71").Value = 4300.6665
$ws.Range("I71").Value = 1826
$ws.Range("J71").Value = 9250
$ws.Range("K71").Value = 9130
$ws.Range("L71").Value = 46250
$ws.Range("M71").Value = -5386
$ws.Range("N71").Value = -53738

$ws.Range("H100").Value = 5055.625
$ws.Range("I100").Value = 6689
$ws.Range("J100").Value = 2333.3333
$ws.Range("K100").Value = 6689
$ws.Range("L100").Value = 2333.3333
$ws.Range("M100").Value = -6148
$ws.Range("N100").Value = -3415.3333

$ws.Range("H126").Value = 50001.773
$ws.Range("J126").Value = 95889.73
$ws.Range("L126").Value = 287669.19
$ws.Range("N126").Value = -292609.19

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 123571.43
$ws.Range("J2").Value = 37501
$ws.Range("L2").Value = 37501
$ws.Range("N2").Value = -37725

$ws.Range("H29").Value = 75000
$ws.Range("J29").Value = 75000
$ws.Range("L29").Value = 75000
$ws.Range("N29").Value = -75580

$ws.Range("H107").Value = 41668690
